$wb = $excel.ActiveWorkbook

# --- Final_Total sheet ---
$ws = $wb.Worksheets.Item("Final_Total")
$ws.Range("D64").Value = 21683.9361904972
$ws.Range("D65").Value = 125572.693809503
$ws.Range("D66").Value = 19068.7296387063
$ws.Range("D67").Value = 79753.4903612937
$ws.Range("D68").Value = 105386.910561649
$ws.Range("D69").Value = 747650.279629641
$ws.Range("D70").Value = 196577.75980871
$ws.Range("D71").Value = 13227.6073808706
$ws.Range("D72").Value = 127121.182619129

# --- Final_Columbus sheet ---
$ws = $wb.Worksheets.Item("Final_Columbus")
$ws.Range("D20").Value = 21683.9361904972
$ws.Range("D21").Value = 125572.693809503

# --- Final_Gahanna sheet ---
$ws = $wb.Worksheets.Item("Final_Gahanna")
$ws.Range("D20").Value = 19068.7296387063
$ws.Range("D21").Value = 79753.4903612937

# --- Final_JeffersonUnincorporated sheet ---
$ws = $wb.Worksheets.Item("Final_JeffersonUnincorporated")
$ws.Range("D20").Value = 105386.910561649
$ws.Range("D21").Value = 747650.279629641
$ws.Range("D22").Value = 196577.75980871

# --- Final_Reynoldsburg sheet ---
$ws = $wb.Worksheets.Item("Final_Reynoldsburg")
$ws.Range("D10").Value = 13227.6073808706
$ws.Range("D11").Value = 127121.182619129
